$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13 (which only held the "1720367 - Teresa Cristina Brazil de Paiva" value
# in B/C, with blank A) is removed entirely; everything below shifts up.
$ws.Rows("13").Delete()

# Row 10 (Objetivos:) now shows the lecturer name instead of the old ementa text.
$ws.Range("B10").Value = "1720367 - Teresa Cristina Brazil de Paiva"
$ws.Range("C10").Value = "1720367 - Teresa Cristina Brazil de Paiva"

# Row 13 (Programa resumido:) now reads "Semestral".
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# Row 15 (Programa:) now reads the activation date.
$ws.Range("B15").Value = "01/01/2018"
$ws.Range("C15").Value = "01/01/2018"

# Row 18 (Método:) now shows the lecturer name.
$ws.Range("B18").Value = "1720367 - Teresa Cristina Brazil de Paiva"
$ws.Range("C18").Value = "1720367 - Teresa Cristina Brazil de Paiva"

# Row 19 (Critério:) now shows the old "Método" evaluation text.
$ws.Range("B19").Value = "Os alunos serão avaliados por meio de duas provas (P1 e P2) e complementada por meio de trabalhos, seminários e/ou relatórios (C)."
$ws.Range("C19").Value = "Os alunos serão avaliados por meio de duas provas (P1 e P2) e complementada por meio de trabalhos, seminários e/ou relatórios (C)."

# Row 20 (Norma de recuperação:) now shows the old "Critério" text.
$ws.Range("B20").Value = "A nota final (NF) será calculada atribuindo-se peso um para a primeira avaliação (P1 = 7 pontos e C = 3 pontos) e peso dois para a segunda avaliação (P2 = 10 pontos).A média ponderada das notas corresponderá à média do período letivo, ou seja: Média do período letivo normal = ((P1 + C) + P2.2)/3.Serão aprovados os alunos que obtiverem média igual ou maior que 5,0 e 70% de frequência no curso."
$ws.Range("C20").Value = "A nota final (NF) será calculada atribuindo-se peso um para a primeira avaliação (P1 = 7 pontos e C = 3 pontos) e peso dois para a segunda avaliação (P2 = 10 pontos).A média ponderada das notas corresponderá à média do período letivo, ou seja: Média do período letivo normal = ((P1 + C) + P2.2)/3.Serão aprovados os alunos que obtiverem média igual ou maior que 5,0 e 70% de frequência no curso."

# Row 21 (Bibliografia:) now shows the old "Norma de recuperação" text.
$ws.Range("B21").Value = "Aos alunos que obtiverem média igual ou maior que 3,0 e menor que 5,0 será oferecido um programa de recuperação que será avaliado por uma prova final. Nesse caso, a média final do aluno será: Média final = (média do período letivo normal + nota prova final)/2.Serão aprovados os alunos que obtiverem média final igual ou maior que 5,0."
$ws.Range("C21").Value = "Aos alunos que obtiverem média igual ou maior que 3,0 e menor que 5,0 será oferecido um programa de recuperação que será avaliado por uma prova final. Nesse caso, a média final do aluno será: Média final = (média do período letivo normal + nota prova final)/2.Serão aprovados os alunos que obtiverem média final igual ou maior que 5,0."
